$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended at the bottom of the tracking list (rows 241-252).
# styleRow indicates which existing formatted row (237 -> "A" style 7,
# 238 -> "A" style 8) to copy the number-format/border/fill from for column A;
# column B/C formatting is copied from row 240 (style 6 / style 11), which is
# constant across all of the new rows.
$newRows = @(
    @(101283, "래미안하이어스", 84, 237),
    @(2895,   "솔거대림", 90, 237),
    @(3864,   "세종주공6단지", 84, 238),
    @(121277, "힐스테이트금정역(주상복합)", 84, 237),
    @(3580,   "금정쌍용", 59, 238),
    @(8386,   "율곡주공3단지", 59, 238),
    @(26398,  "군포대야미e-편한세상", 59, 237),
    @(2886,   "가야주공5단지1차", 58, 238),
    @(7963,   "한양수리", 84, 237),
    @(101480, "의왕역센트럴시티", 84, 237),
    @(104999, "청천마을대우", 84, 238),
    @(8333,   "무지개마을대림", 84, 238)
)

$startRow = 241
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $styleSrcRow = $row[3]

    # Copy formatting first (number format, borders, fill, font) from the
    # matching template rows, column by column.
    $ws.Range("A$styleSrcRow").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("B240").Copy()
    $ws.Range("B$r").PasteSpecial(-4122)
    $ws.Range("C240").Copy()
    $ws.Range("C$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    $ws.Rows.Item($r).RowHeight = 13.8
}

$excel.CutCopyMode = 0

$ws.Range("A226").Select()
$ws.Range("M251").Select()
